# Update column F (dSF) values to reflect repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -3
$ws.Range("F3").Value  = -3
$ws.Range("F5").Value  = -24
$ws.Range("F6").Value  = -3
$ws.Range("F7").Value  = -7
$ws.Range("F8").Value  = 1
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = -4
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = 3
